$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting rows 25:33 down to 26:34
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with values, matching the pattern of adjacent rows
$ws.Cells.Item(25, 1).Value = 5
$ws.Cells.Item(25, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(25, 3).Value = "Maule"
$ws.Cells.Item(25, 4).Value = 44489
$ws.Cells.Item(25, 5).Value = 7
$ws.Cells.Item(25, 6).Value = 300000000
$ws.Cells.Item(25, 7).Value = "Espárragos"
$ws.Cells.Item(25, 8).Value = "Verde"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 4000
$ws.Cells.Item(25, 11).Value = 900
$ws.Cells.Item(25, 12).Value = 900
$ws.Cells.Item(25, 13).Value = 900
$ws.Cells.Item(25, 14).Value = "$/kilo"
$ws.Cells.Item(25, 15).Value = "Provincia de Linares"
$ws.Cells.Item(25, 16).Value = 900
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Apply the same date number format style as other D-column cells
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
